# ChemicalDatabase.xlsx - "Add files via upload" edit
#
# The author trimmed the chemical-properties table on Sheet1: the rows for
# Water, Carbon dioxide, and the various alkali/alkaline-earth salts
# (rows 51-65, shared-string entries "Water" ... "Ca(NO3)2") were removed,
# leaving two empty-but-styled rows (51-52) in their place and shrinking
# the used range down to A1:E52. A handful of rows with long, two-line
# chemical names (Potassium hydroxide, Calcium carbonate, Sodium
# bicarbonate, Sodium Hypochlorite, Magnesium sulfate) ended up with a
# taller row height, and the view was left scrolled near the top with
# zoom back at 100% and a different cell selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose wrapped two-line chemical name now needs the taller 30pt row.
$ws.Rows.Item(14).RowHeight = 30   # Potassium hydroxide
$ws.Rows.Item(38).RowHeight = 30   # Calcium carbonate
$ws.Rows.Item(42).RowHeight = 30   # Sodium bicarbonate
$ws.Rows.Item(43).RowHeight = 30   # Sodium Hypochlorite
$ws.Rows.Item(46).RowHeight = 30   # Magnesium sulfate

# Drop the trailing chemical entries (Water through Calcium nitrate).
# Clearing A51:E65 first leaves rows 51-52 present but blank (matching the
# target, which keeps two styled-but-empty rows), then the now-unused
# rows 53-65 are removed outright so the sheet ends at row 52.
$ws.Range("A51:E65").ClearContents()
$ws.Rows("53:65").Delete()

# Leave the view the way the saved workbook shows it: zoomed to 100%
# (back down from the 102% it was saved at) with a different cell active.
$excel.ActiveWindow.Zoom = 100
$ws.Range("O48").Select()
